$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update simple in-place values on the "Metadata" sheet ---
$ws1.Range("B3").Value2  = "0.1.7"
$ws1.Range("B6").Value2  = "draft"
$ws1.Range("B8").Value2  = "2024-11-22T12:33:30-06:00"
$ws1.Range("B10").Value2 = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws1.Range("B11").Value2 = "Bob Milius (bmilius@nmdp.org)"

# --- Make room for the new "Jurisdiction" row at position 12 by shifting the
#     existing rows 12-15 (Description, Purpose, Copyright, Immutable) down
#     by one, to rows 13-16. Work bottom-up so we never clobber a value
#     before it has been copied down. Formatting (style 2) is already applied
#     to row 16 from the original A1:B15 range via the sheet defaults, but we
#     also paste formats explicitly to be safe. ---

$ws1.Range("A16").Value2 = $ws1.Range("A15").Value2
$ws1.Range("B16").Value2 = $ws1.Range("B15").Value2

$ws1.Range("A15").Value2 = $ws1.Range("A14").Value2
$ws1.Range("B15").Value2 = $ws1.Range("B14").Value2

$ws1.Range("A14").Value2 = $ws1.Range("A13").Value2
$ws1.Range("B14").Value2 = $ws1.Range("B13").Value2

$ws1.Range("A13").Value2 = $ws1.Range("A12").Value2
$ws1.Range("B13").Value2 = $ws1.Range("B12").Value2

# --- Row 12 now becomes the new "Jurisdiction" row (empty value) ---
$ws1.Range("A12").Value2 = "Jurisdiction"
$ws1.Range("B12").Value2 = ""

# --- Make sure every cell in the table (A1:B16) carries the same formatting
#     style used throughout the sheet (header row 1 uses style 1, all other
#     rows use style 2). Copy formats from row 15 (style 2) down onto row 16
#     explicitly so the newly extended row picks up the same look as the
#     rest of the table. ---
$ws1.Range("A15:B15").Copy() | Out-Null
$ws1.Range("A16:B16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
